# Update "想去人数" (column F) figures across all four sheets to match
# newly generated output (gh-pages commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 1396
$ws.Cells.Item(5, 6).Value = 5798
$ws.Cells.Item(6, 6).Value = 485
$ws.Cells.Item(8, 6).Value = 7
$ws.Cells.Item(9, 6).Value = 3433
$ws.Cells.Item(10, 6).Value = 6657
$ws.Cells.Item(11, 6).Value = 208
$ws.Cells.Item(12, 6).Value = 1316
$ws.Cells.Item(13, 6).Value = 762
$ws.Cells.Item(14, 6).Value = 97
$ws.Cells.Item(17, 6).Value = 1116
$ws.Cells.Item(19, 6).Value = 110
$ws.Cells.Item(21, 6).Value = 172
$ws.Cells.Item(23, 6).Value = 974
$ws.Cells.Item(24, 6).Value = 322
$ws.Cells.Item(27, 6).Value = 108
$ws.Cells.Item(30, 6).Value = 19
$ws.Cells.Item(31, 6).Value = 52
$ws.Cells.Item(32, 6).Value = 2
$ws.Cells.Item(35, 6).Value = 301
$ws.Cells.Item(36, 6).Value = 13
$ws.Cells.Item(38, 6).Value = 291
$ws.Cells.Item(39, 6).Value = 1167
$ws.Cells.Item(40, 6).Value = 52
$ws.Cells.Item(41, 6).Value = 97

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(6, 6).Value = 513
$ws.Cells.Item(11, 6).Value = 114
$ws.Cells.Item(12, 6).Value = 9
$ws.Cells.Item(16, 6).Value = 1721
$ws.Cells.Item(21, 6).Value = 183
$ws.Cells.Item(23, 6).Value = 133
$ws.Cells.Item(26, 6).Value = 607
$ws.Cells.Item(30, 6).Value = 682
$ws.Cells.Item(31, 6).Value = 954
$ws.Cells.Item(32, 6).Value = 568
$ws.Cells.Item(34, 6).Value = 82
$ws.Cells.Item(37, 6).Value = 95
$ws.Cells.Item(38, 6).Value = 123

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 701
$ws.Cells.Item(5, 6).Value = 825
$ws.Cells.Item(6, 6).Value = 548
$ws.Cells.Item(7, 6).Value = 280
$ws.Cells.Item(8, 6).Value = 1042

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 701
$ws.Cells.Item(4, 6).Value = 1396
$ws.Cells.Item(5, 6).Value = 825
$ws.Cells.Item(8, 6).Value = 548
$ws.Cells.Item(9, 6).Value = 548
$ws.Cells.Item(10, 6).Value = 280
$ws.Cells.Item(11, 6).Value = 280
$ws.Cells.Item(12, 6).Value = 513
$ws.Cells.Item(14, 6).Value = 5798
$ws.Cells.Item(15, 6).Value = 485
$ws.Cells.Item(17, 6).Value = 3433
$ws.Cells.Item(19, 6).Value = 6658
$ws.Cells.Item(20, 6).Value = 208
$ws.Cells.Item(21, 6).Value = 1316
$ws.Cells.Item(24, 6).Value = 762
$ws.Cells.Item(25, 6).Value = 97
$ws.Cells.Item(26, 6).Value = 1042
$ws.Cells.Item(27, 6).Value = 183
$ws.Cells.Item(29, 6).Value = 1116
$ws.Cells.Item(30, 6).Value = 110
$ws.Cells.Item(31, 6).Value = 172
$ws.Cells.Item(32, 6).Value = 974
$ws.Cells.Item(33, 6).Value = 607
$ws.Cells.Item(34, 6).Value = 322
$ws.Cells.Item(36, 6).Value = 108
$ws.Cells.Item(38, 6).Value = 19
$ws.Cells.Item(39, 6).Value = 52
$ws.Cells.Item(42, 6).Value = 954
$ws.Cells.Item(43, 6).Value = 568
$ws.Cells.Item(44, 6).Value = 301
$ws.Cells.Item(45, 6).Value = 82
$ws.Cells.Item(46, 6).Value = 291
$ws.Cells.Item(47, 6).Value = 95
$ws.Cells.Item(48, 6).Value = 123
$ws.Cells.Item(50, 6).Value = 97

$wb.Save()
